$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; this shifts existing rows 33-75 down to 34-76
# and keeps formatting (e.g. date style on column D) consistent with the diff.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record's data
$ws.Cells.Item(33, 1).Value = 5
$ws.Cells.Item(33, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(33, 3).Value = "Maule"
$ws.Cells.Item(33, 4).Value = 44797
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = 100112026
$ws.Cells.Item(33, 7).Value = "Haba"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 10000
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = 10000
$ws.Cells.Item(33, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(33, 16).Value = 400
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"
